# The deck ships two embedded theme parts:
#   ppt/theme/theme1.xml -> currently the stock "Office Theme" colours
#   ppt/theme/theme2.xml -> currently the "Integral" theme colours, and this
#                            is the theme actually driving the slide master /
#                            every slide in the deck.
# The authored edit swaps the two themes' contents so the deck (theme2.xml)
# now carries the standard "Office" colour palette instead of "Integral".
#
# PowerPoint's automation model doesn't expose a "replace this theme part
# with that one" verb, so we reproduce the net effect the supported way:
# push the 12 Office-theme RGB values into the live ThemeColorScheme, which
# is backed by the clrScheme actually referenced by the slides (theme2.xml).
# Index order follows MsoThemeColorSchemeIndex: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      #000000
$tcs.Item(2).RGB  = 16777215   # lt1      #FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      #44546A
$tcs.Item(4).RGB  = 15132391   # lt2      #E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  #5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  #ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  #A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  #FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  #4472C4
$tcs.Item(10).RGB = 4697456    # accent6  #70AD47
$tcs.Item(11).RGB = 12673797   # hlink    #0563C1
$tcs.Item(12).RGB = 7491477    # folHlink #954F72
